# The deck's design ("Integral") theme (ppt/theme/theme1.xml, the theme
# actually applied to the slide master / all slides) is being swapped for
# the stock "Office Theme" palette (the colours previously only used by
# the notes master's theme part). Font scheme / format scheme are already
# identical between the two themes, so the only externally-visible change
# is the 12-slot theme colour scheme (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink).
#
# PowerPoint's object model doesn't give VBA/COM a "rename theme" or
# "swap theme parts" verb, but ActiveX/VBA automation DOES let you repaint
# every slot of a slide's ThemeColorScheme - that's exactly the gesture a
# user makes from Design > Variants > Colors > Customize Colors, and it
# rewrites the theme part backing the slide master in place.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# New "Office Theme" swatches (RGB packed as PowerPoint's 0x00BBGGRR long):
$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
